$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction -------------------------------------------------
# Row 19 (element_attribute_name) previously held the attribute *value*
# string ("font-size: 18px;font-weight:bold;color:#000000") by mistake;
# it should hold the attribute *name* ("style").
$ws.Range("B19").Value = "style"

# Row 20 (element_attribute_value) was left blank; it should hold the
# value that row 19 used to (incorrectly) contain.
$ws.Range("B20").Value = "font-size: 18px;font-weight:bold;color:#000000"

# --- View / selection state -------------------------------------------
# Move the active selection to C19 and scroll the window so row 9 is at
# the top (matches the saved sheetView in the authored workbook).
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C19").Select()
